# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Update cryptocurrency price and volume(1h) data
Set-TextValue "D2" "25.679.26"
Set-TextValue "E2" "  -3.50%  "

Set-TextValue "D3" "1.744.65"
Set-TextValue "E3" "  -5.57%  "

Set-TextValue "E4" "  +0.01%  "

Set-TextValue "D5" "235.45"
Set-TextValue "E5" "  -10.22%  "

Set-TextValue "E6" "  -0.09%  "

Set-TextValue "D7" "0.4927"
Set-TextValue "E7" "  -7.53%  "

Set-TextValue "D8" "41.60"
Set-TextValue "E8" "  -7.69%  "

Set-TextValue "D9" "0.2540"
Set-TextValue "E9" "  -19.69%  "

Set-TextValue "D10" "0.06011"
Set-TextValue "E10" "  -13.76%  "

Set-TextValue "D11" "1.744.36"
Set-TextValue "E11" "  -5.61%  "

Set-TextValue "D12" "0.06837"
Set-TextValue "E12" "  -12.82%  "

Set-TextValue "D13" "14.75"
Set-TextValue "E13" "  -21.89%  "

Set-TextValue "D14" "4.449"
Set-TextValue "E14" "  -11.77%  "

Set-TextValue "D15" "76.72"
Set-TextValue "E15" "  -14.27%  "

Set-TextValue "D16" "0.5658"
Set-TextValue "E16" "  -26.77%  "

Set-TextValue "D17" "1.001"

Set-TextValue "D19" "25.720.33"
Set-TextValue "E19" "  -3.42%  "

Set-TextValue "D20" "11.27"
Set-TextValue "E20" "  -20.31%  "

Set-TextValue "D21" "0.000006555"
Set-TextValue "E21" "  -17.89%  "

Set-TextValue "D22" "1.966.17"
Set-TextValue "E22" "  -5.73%  "

Set-TextValue "D23" "4.003"
Set-TextValue "E23" "  -13.87%  "

Set-TextValue "D24" "5.013"
Set-TextValue "E24" "  -16.82%  "

Set-TextValue "E25" "  -16.09%  "

Set-TextValue "D26" "136.75"
Set-TextValue "E26" "  -3.88%  "

Set-TextValue "D27" "1.482"
Set-TextValue "E27" "  -12.93%  "

Set-TextValue "D28" "1.815"
Set-TextValue "E28" "  -18.03%  "

Set-TextValue "D29" "14.65"
Set-TextValue "E29" "  -14.56%  "

Set-TextValue "D30" "101.90"
Set-TextValue "E30" "  -8.82%  "

Set-TextValue "D31" "3.755"
Set-TextValue "E31" "  -13.09%  "

Set-TextValue "D32" "0.07974"
Set-TextValue "E32" "  -9.18%  "

Set-TextValue "D33" "3.370"
Set-TextValue "E33" "  -18.09%  "

Set-TextValue "D34" "0.04379"
Set-TextValue "E34" "  -9.99%  "

Set-TextValue "D35" "1.000"
Set-TextValue "E35" "  -0.09%  "

Set-TextValue "D36" "2.616"
Set-TextValue "E36" "  -9.41%  "

Set-TextValue "D37" "0.9735"
Set-TextValue "E37" "  -14.59%  "

Set-TextValue "D38" "0.6020"
Set-TextValue "E38" "  -18.70%  "

Set-TextValue "D39" "2.657"
Set-TextValue "E39" "  -14.41%  "

Set-TextValue "D40" "2.004"
Set-TextValue "E40" "  -15.00%  "

Set-TextValue "E41" "  -0.04%  "

Set-TextValue "D42" "102.62"
Set-TextValue "E42" "  -5.78%  "

Set-TextValue "D43" "0.01500"
Set-TextValue "E43" "  -13.78%  "

Set-TextValue "D44" "0.7530"
Set-TextValue "E44" "  -16.97%  "

Set-TextValue "D45" "5.157"
Set-TextValue "E45" "  -12.78%  "

Set-TextValue "D46" "0.3718"
Set-TextValue "E46" "  -22.96%  "

Set-TextValue "D47" "0.05235"
Set-TextValue "E47" "  -10.02%  "

Set-TextValue "D48" "0.1063"
Set-TextValue "E48" "  -15.02%  "

Set-TextValue "D49" "29.97"
Set-TextValue "E49" "  -14.61%  "

Set-TextValue "D50" "5.879"
Set-TextValue "E50" "  -23.79%  "

Set-TextValue "D51" "52.22"
Set-TextValue "E51" "  -13.70%  "
